# Applies the edits described by the commit "Archivo de selección subido"
# to the "Seleccion.xlsx" workbook (sheet "Hoja1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update the selected robot candidate (E21): IRB 120 -> IRB 1200 ---
$ws.Range("E21").Value = "IRB 1200"

# --- Split the "Alcance" value (D23) into a number + unit ---
# Previously D23 held the text "0.58 m"; now the magnitude goes in D23
# as a real number and the unit "m" moves to the new cell E23.
$ws.Range("D23").Value = 0.9
$ws.Range("E23").Value = "m"

# --- Update the view state (zoom level, scroll position, selection) ---
$win = $excel.ActiveWindow
$win.Zoom = 115
$win.ScrollRow = 20
$win.ScrollColumn = 1
[void]$ws.Range("E24").Select()
